# Tadata/funder_data.xlsx - "Add files via upload"
#
# The funder account row for 001302843 / FP0051 (USD) was removed from the
# funder list on Sheet1. Deleting the whole row shifts every row below it
# up by one, which is exactly what the target workbook shows (rows that
# used to be 16-25 are now 15-24), and Excel automatically keeps the
# shared-strings table, cell styles and the sheet dimension in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 holds account 001302843 / FP0051 - remove it entirely.
$ws.Rows(15).Delete() | Out-Null

# Leave the selection where the author last clicked after the edit.
$ws.Range("D20").Select() | Out-Null
